$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, using the same style as the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-28
$data = @{
    2  = @(5, 6)
    3  = @(6, 6)
    4  = @(8, 8)
    5  = @(1, 1)
    6  = @(1, 3)
    7  = @(7, 7)
    8  = @(5, 6)
    9  = @(7, 7)
    10 = @(5, 7)
    11 = @(5, 5)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(5, 5)
    15 = @(8, 8)
    16 = @(4, 5)
    17 = @(7, 7)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(8, 9)
    21 = @(5, 5)
    22 = @(4, 4)
    23 = @(7, 7)
    24 = @(5, 5)
    25 = @(5, 5)
    26 = @(7, 7)
    27 = @(7, 7)
    28 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

$wb.Save()
